$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of Kiwi price data (2022-07-27, serial 44769) was added to the
# "logica_diaria" dataset. In this sheet it lands as 5 rows (one per quality
# grade) inserted right before the existing row 814, pushing all subsequent
# rows down by 5 (old 814 -> new 819, ..., old 901 -> new 906) and extending
# the used range from A1:T901 to A1:T906.
$ws.Range("A814:A818").EntireRow.Insert()

$ws.Range("A814:T814").Value = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44769, 13, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Especial", 20, 200000, 200000, 200000, "$/bins (450 kilos)", "Región de O'Higgins", 444, 450)

$ws.Range("A815:T815").Value = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44769, 13, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Extra (doble especial)", 18, 300000, 300000, 300000, "$/bins (450 kilos)", "Región de O'Higgins", 667, 450)

$ws.Range("A816:T816").Value = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44769, 13, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Primera", 38, 160000, 170000, 165000, "$/bins (450 kilos)", "Región de O'Higgins", 367, 450)

$ws.Range("A817:T817").Value = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44769, 13, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Segunda", 25, 120000, 120000, 120000, "$/bins (450 kilos)", "Región de O'Higgins", 267, 450)

$ws.Range("A818:T818").Value = @(6, "Mercado Mayorista Lo Valledor de Santiago", "Metropolitana", 44769, 13, "Fruta", 100101, "Berries", 100101007, "Kiwi", "Hayward", "Tercera", 20, 100000, 100000, 100000, "$/bins (450 kilos)", "Región de O'Higgins", 222, 450)
